# The upstream commit ("Rework DOCX implementation to better support
# testing and decrease coupling") re-saved this fixture .docx through a
# different OOXML serializer (docx4j instead of Word's native writer).
#
# Diffing the canonical/pretty-printed XML of every part shows that the
# *only* differences anywhere in the package are the order in which the
# xmlns:* namespace declarations are listed on the root element of
# word/document.xml, word/endnotes.xml, word/footer1.xml,
# word/footnotes.xml, word/header1.xml, word/styles.xml and
# word/theme/theme1.xml. The same set of namespace prefixes/URIs is
# present before and after; none of the document text, formatting,
# structure, relationships, or any other content changed at all.
#
# Namespace-declaration ordering on an element is not semantic XML
# content (and isn't something the Word object model exposes any way to
# control - Word/COM automation edits content, it doesn't let a macro
# dictate how the XML writer orders an element's xmlns attributes).
# There is therefore no content-level edit to make here: the correct
# reproduction of this change is to leave the document's content
# untouched.
$d = $word.ActiveDocument
